# Scheduled-runner price/profit refresh across the Leve-profit sheets.
# Updates market-price-derived columns (H:N) for the rows whose
# currentAveragePrice/NQ/HQ and downstream profit figures moved since the
# last pull. A couple of rows gain/lose a trailing profit cell entirely
# because the NQ/HQ price split changed (e.g. ALC!M125 no longer applies,
# CUL!N97 newly does).

$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 3756.2354
$ws.Range("I28").Value = 2036
$ws.Range("J28").Value = 5691.5
$ws.Range("K28").Value = 2036
$ws.Range("L28").Value = 5691.5
$ws.Range("M28").Value = -1551
$ws.Range("N28").Value = -6661.5
$ws.Range("H33").Value = 272.27274
$ws.Range("J33").Value = 99.666664
$ws.Range("L33").Value = 99.666664
$ws.Range("N33").Value = -557.666664
$ws.Range("H41").Value = 1297.421
$ws.Range("I41").Value = 700.1429000000001
$ws.Range("J41").Value = 1645.8334
$ws.Range("K41").Value = 700.1429000000001
$ws.Range("L41").Value = 1645.8334
$ws.Range("M41").Value = -260.1429000000001
$ws.Range("N41").Value = -2525.8334
$ws.Range("H116").Value = 3740097
$ws.Range("I116").Value = 5162311
$ws.Range("J116").Value = 6785.25
$ws.Range("K116").Value = 5162311
$ws.Range("L116").Value = 6785.25
$ws.Range("M116").Value = -5158869
$ws.Range("N116").Value = -13669.25
$ws.Range("H121").Value = 4243.143
$ws.Range("J121").Value = 4243.143
$ws.Range("L121").Value = 12729.429
$ws.Range("N121").Value = -16223.429
$ws.Range("H125").Value = 4187.2
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 4187.2
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 37684.8
$ws.Range("M125").Value = ""
$ws.Range("N125").Value = -42604.8
$ws.Range("H135").Value = 2459.0557
$ws.Range("I135").Value = 2766.182
$ws.Range("J135").Value = 1976.4286
$ws.Range("K135").Value = 24895.638
$ws.Range("L135").Value = 17787.8574
$ws.Range("M135").Value = -22360.638
$ws.Range("N135").Value = -22857.8574
$ws.Range("H137").Value = 3098.4285
$ws.Range("I137").Value = 2354.625
$ws.Range("J137").Value = 3318.8147
$ws.Range("K137").Value = 7063.875
$ws.Range("L137").Value = 9956.444100000001
$ws.Range("M137").Value = -4513.875
$ws.Range("N137").Value = -15056.4441
$ws.Range("H138").Value = 4987.353
$ws.Range("I138").Value = 2458.2632
$ws.Range("J138").Value = 5968.0205
$ws.Range("K138").Value = 7374.7896
$ws.Range("L138").Value = 17904.0615
$ws.Range("M138").Value = -2234.7896
$ws.Range("N138").Value = -28184.0615

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3241.348
$ws.Range("I32").Value = 2644.7144
$ws.Range("K32").Value = 2644.7144
$ws.Range("M32").Value = -2357.7144
$ws.Range("H41").Value = 13688.333
$ws.Range("I41").Value = 4032.75
$ws.Range("K41").Value = 4032.75
$ws.Range("M41").Value = -3618.75
$ws.Range("H61").Value = 13093.979
$ws.Range("I61").Value = 13007.15
$ws.Range("K61").Value = 13007.15
$ws.Range("M61").Value = -12795.15
$ws.Range("H74").Value = 5817108.5
$ws.Range("I74").Value = 13159365
$ws.Range("J74").Value = 4488.9165
$ws.Range("K74").Value = 13159365
$ws.Range("L74").Value = 4488.9165
$ws.Range("M74").Value = -13158491
$ws.Range("N74").Value = -6236.9165
$ws.Range("H77").Value = 5817108.5
$ws.Range("I77").Value = 13159365
$ws.Range("J77").Value = 4488.9165
$ws.Range("K77").Value = 65796825
$ws.Range("L77").Value = 22444.5825
$ws.Range("M77").Value = -65792457
$ws.Range("N77").Value = -31180.5825
$ws.Range("H132").Value = 14088.985
$ws.Range("I132").Value = 18604.39
$ws.Range("K132").Value = 55813.17
$ws.Range("M132").Value = -53283.17
$ws.Range("H133").Value = 83680.5
$ws.Range("J133").Value = 83680.5
$ws.Range("L133").Value = 83680.5
$ws.Range("N133").Value = -88740.5
$ws.Range("H136").Value = 13093.979
$ws.Range("I136").Value = 13007.15
$ws.Range("K136").Value = 39021.45
$ws.Range("M136").Value = -36471.45

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2117.35
$ws.Range("I86").Value = 2117.2
$ws.Range("J86").Value = 2117.8
$ws.Range("K86").Value = 2117.2
$ws.Range("L86").Value = 2117.8
$ws.Range("M86").Value = -994.1999999999998
$ws.Range("N86").Value = -4363.8
$ws.Range("H89").Value = 2117.35
$ws.Range("I89").Value = 2117.2
$ws.Range("J89").Value = 2117.8
$ws.Range("K89").Value = 10586
$ws.Range("L89").Value = 10589
$ws.Range("M89").Value = -4970
$ws.Range("N89").Value = -21821
$ws.Range("H107").Value = 1867.3928
$ws.Range("I107").Value = 1713.8948
$ws.Range("J107").Value = 2191.4443
$ws.Range("K107").Value = 1713.8948
$ws.Range("L107").Value = 2191.4443
$ws.Range("M107").Value = 206.1052
$ws.Range("N107").Value = -6031.4443

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20836608
$ws.Range("I31").Value = 41668204
$ws.Range("J31").Value = 5012.875
$ws.Range("K31").Value = 41668204
$ws.Range("L31").Value = 5012.875
$ws.Range("M31").Value = -41667909
$ws.Range("N31").Value = -5602.875
$ws.Range("H34").Value = 20836608
$ws.Range("I34").Value = 41668204
$ws.Range("J34").Value = 5012.875
$ws.Range("K34").Value = 41668204
$ws.Range("L34").Value = 5012.875
$ws.Range("M34").Value = -41668002
$ws.Range("N34").Value = -5416.875
$ws.Range("H81").Value = 74600
$ws.Range("J81").Value = 74600
$ws.Range("L81").Value = 74600
$ws.Range("N81").Value = -76596
$ws.Range("H84").Value = 74600
$ws.Range("J84").Value = 74600
$ws.Range("L84").Value = 223800
$ws.Range("N84").Value = -233784

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1124.9565
$ws.Range("I5").Value = 895
$ws.Range("J5").Value = 1335.75
$ws.Range("K5").Value = 2685
$ws.Range("L5").Value = 4007.25
$ws.Range("M5").Value = -2573
$ws.Range("N5").Value = -4231.25
$ws.Range("H97").Value = 323.25
$ws.Range("I97").Value = 401.5
$ws.Range("J97").Value = 245
$ws.Range("K97").Value = 1204.5
$ws.Range("L97").Value = 735
$ws.Range("M97").Value = -708.5
$ws.Range("N97").Value = -1727
$ws.Range("H135").Value = 1124.9565
$ws.Range("I135").Value = 895
$ws.Range("J135").Value = 1335.75
$ws.Range("K135").Value = 8055
$ws.Range("L135").Value = 12021.75
$ws.Range("M135").Value = -5520
$ws.Range("N135").Value = -17091.75

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1286654.9
$ws.Range("I80").Value = 1851837.8
$ws.Range("J80").Value = 14993.5
$ws.Range("K80").Value = 1851837.8
$ws.Range("L80").Value = 14993.5
$ws.Range("M80").Value = -1850839.8
$ws.Range("N80").Value = -16989.5
$ws.Range("H83").Value = 1286654.9
$ws.Range("I83").Value = 1851837.8
$ws.Range("J83").Value = 14993.5
$ws.Range("K83").Value = 9259189
$ws.Range("L83").Value = 74967.5
$ws.Range("M83").Value = -9254197
$ws.Range("N83").Value = -84951.5
$ws.Range("H126").Value = 3295.3713
$ws.Range("I126").Value = 2144.25
$ws.Range("J126").Value = 5806.909
$ws.Range("K126").Value = 6432.75
$ws.Range("L126").Value = 17420.727
$ws.Range("M126").Value = -3962.75
$ws.Range("N126").Value = -22360.727
$ws.Range("H132").Value = 3138.547
$ws.Range("I132").Value = 2776.814
$ws.Range("K132").Value = 8330.441999999999
$ws.Range("M132").Value = -5800.441999999999
$ws.Range("H133").Value = 119274.5
$ws.Range("J133").Value = 119274.5
$ws.Range("L133").Value = 119274.5
$ws.Range("N133").Value = -129394.5
$ws.Range("H139").Value = 90000
$ws.Range("I139").Value = 90000
$ws.Range("K139").Value = 90000
$ws.Range("M139").Value = -84860

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5210766
$ws.Range("J82").Value = 5601.5
$ws.Range("L82").Value = 5601.5
$ws.Range("N82").Value = -6323.5
$ws.Range("H85").Value = 5210766
$ws.Range("J85").Value = 5601.5
$ws.Range("L85").Value = 5601.5
$ws.Range("N85").Value = -8097.5
$ws.Range("H132").Value = 3195.31
$ws.Range("I132").Value = 3195.1753
$ws.Range("J132").Value = 3199.6667
$ws.Range("K132").Value = 9585.525900000001
$ws.Range("L132").Value = 9599.000100000001
$ws.Range("M132").Value = -7055.525900000001
$ws.Range("N132").Value = -14659.0001
$ws.Range("H136").Value = 3950.39
$ws.Range("I136").Value = 3965.0403
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 11895.1209
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -9345.1209
$ws.Range("N136").Value = -12600

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 70000
$ws.Range("J94").Value = 70000
$ws.Range("L94").Value = 70000
$ws.Range("N94").Value = -71802
